# Updates cryptos list data (price + 1h volume change figures) and
# applies two row reshuffles that occurred in the upstream scrape:
#   - rows 17/18 swap (Uniswap <-> WrappedEther)
#   - rows 47-50 rotate (ThetaToken/Celestia/EnergySwap/PEPE)
# A leading "'" forces numeric-looking text (e.g. "0.730") to stay text,
# matching the sheet's existing inline-string cells; the follow-up
# Style = "Normal" strips the quote-prefix formatting flag Excel would
# otherwise leave behind, keeping cells on the default (unstyled) xf.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Address, [string]$Text)
    $cell = $ws.Range($Address)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

Set-CellText 'D2' '62.299.08'
Set-CellText 'E2' '  +0.50%  '
Set-CellText 'D3' '3.437.85'
Set-CellText 'E3' '  +0.59%  '
Set-CellText 'D4' '0.997'
Set-CellText 'E4' '  -0.36%  '
Set-CellText 'D5' '415.12'
Set-CellText 'E5' '  +1.20%  '
Set-CellText 'D6' '128.35'
Set-CellText 'E6' '  -0.97%  '
Set-CellText 'D7' '0.628'
Set-CellText 'E7' '  -3.25%  '
Set-CellText 'D8' '0.999'
Set-CellText 'E8' '  +0.01%  '
Set-CellText 'D9' '0.730'
Set-CellText 'E9' '  -1.75%  '
Set-CellText 'D10' '0.142'
Set-CellText 'E10' '  +0.17%  '
Set-CellText 'D11' '42.86'
Set-CellText 'E11' '  +0.52%  '
Set-CellText 'E12' '  +2.74%  '
Set-CellText 'D13' '9.22'
Set-CellText 'E13' '  +1.15%  '
Set-CellText 'D14' '3.963.62'
Set-CellText 'E14' '  +0.12%  '
Set-CellText 'E15' '  -0.14%  '
Set-CellText 'D16' '20.63'
Set-CellText 'E16' '  -1.25%  '
Set-CellText 'B17' 'WrappedEther'
Set-CellText 'C17' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText 'D17' '3.463.05'
Set-CellText 'E17' '  +1.61%  '
Set-CellText 'B18' 'Uniswap'
Set-CellText 'C18' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText 'D18' '12.85'
Set-CellText 'E18' '  +6.14%  '
Set-CellText 'D19' '1.08'
Set-CellText 'E19' '  +0.73%  '
Set-CellText 'D20' '62.271.24'
Set-CellText 'E20' '  +0.44%  '
Set-CellText 'D21' '487.53'
Set-CellText 'E21' '  +9.90%  '
Set-CellText 'D22' '91.84'
Set-CellText 'E22' '  +0.38%  '
Set-CellText 'D23' '3.27'
Set-CellText 'E23' '  +3.21%  '
Set-CellText 'D24' '13.18'
Set-CellText 'E24' '  +0.80%  '
Set-CellText 'D25' '3.32'
Set-CellText 'E25' '  +2.05%  '
Set-CellText 'D26' '9.81'
Set-CellText 'E26' '  +11.58%  '
Set-CellText 'D27' '33.65'
Set-CellText 'E27' '  -2.49%  '
Set-CellText 'D28' '4.77'
Set-CellText 'E28' '  +0.33%  '
Set-CellText 'D29' '7.66'
Set-CellText 'E29' '  +0.61%  '
Set-CellText 'D30' '11.88'
Set-CellText 'E30' '  -1.51%  '
Set-CellText 'D31' '2.65'
Set-CellText 'E31' '  -1.40%  '
Set-CellText 'D32' '0.168'
Set-CellText 'E32' '  -1.07%  '
Set-CellText 'E33' '  -2.48%  '
Set-CellText 'D34' '41.03'
Set-CellText 'E34' '  -4.02%  '
Set-CellText 'E35' '  +0.10%  '
Set-CellText 'D36' '58.97'
Set-CellText 'E36' '  +9.52%  '
Set-CellText 'D37' '0.0489'
Set-CellText 'E37' '  -2.73%  '
Set-CellText 'D38' '0.998'
Set-CellText 'E38' '  -0.06%  '
Set-CellText 'D39' '3.05'
Set-CellText 'E39' '  +4.48%  '
Set-CellText 'D40' '0.326'
Set-CellText 'E40' '  +3.88%  '
Set-CellText 'D41' '149.17'
Set-CellText 'E41' '  +5.69%  '
Set-CellText 'D42' '0.135'
Set-CellText 'E42' '  -1.57%  '
Set-CellText 'D43' '3.34'
Set-CellText 'E43' '  -1.01%  '
Set-CellText 'D44' '2.11'
Set-CellText 'E44' '  +6.67%  '
Set-CellText 'D45' '2.60'
Set-CellText 'E45' '  +8.24%  '
Set-CellText 'D46' '4.25'
Set-CellText 'E46' '  +2.69%  '
Set-CellText 'B47' 'PEPE'
Set-CellText 'C47' 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-CellText 'D47' '0.0₃0561'
Set-CellText 'E47' '  +30.70%  '
Set-CellText 'B48' 'ThetaToken'
Set-CellText 'C48' 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-CellText 'D48' '2.33'
Set-CellText 'E48' '  +17.78%  '
Set-CellText 'B49' 'Celestia'
Set-CellText 'C49' 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-CellText 'D49' '16.41'
Set-CellText 'E49' '  -1.39%  '
Set-CellText 'B50' 'EnergySwap'
Set-CellText 'C50' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText 'D50' '22.37'
Set-CellText 'E50' '  -0.03%  '
Set-CellText 'D51' '115.33'
Set-CellText 'E51' '  +7.28%  '
